# CEDS_Code_Change_Tracker.xlsx - "Update Code Change Tracker"
# Adds three new change-log entries (rows 95-97) to Sheet1, matching the
# existing table's layout/formatting, and updates the sheet's frozen-pane
# scroll position / selection to reflect the new bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Date number format used throughout column E/F/G ("Proposed" / "Code Review
# Completed" / "Committed") for the existing rows.
$dateFmt = "mm-dd-yy"

# --- Row 95 ---------------------------------------------------------------
$ws.Range("A95").Value = "Mod A - Fixes energy as driver data"
$ws.Range("B95").Value = 91
$ws.Range("C95").Value = "Rachel Hoesly"
$ws.Range("D95").Value = "Committed"
$ws.Range("E95").Value = 42342
$ws.Range("E95").NumberFormat = $dateFmt
$ws.Range("G95").Value = 42342
$ws.Range("G95").NumberFormat = $dateFmt
$ws.Range("H95").Value = "b0a6d6"

# --- Row 96 ---------------------------------------------------------------
$ws.Range("A96").Value = "Updates addtoDB function, makes faster "
$ws.Range("B96").Value = 92
$ws.Range("C96").Value = "Rachel Hoesly"
$ws.Range("D96").Value = "Committed"
$ws.Range("E96").Value = 42345
$ws.Range("E96").NumberFormat = $dateFmt
$ws.Range("G96").Value = 42345
$ws.Range("G96").NumberFormat = $dateFmt
$ws.Range("H96").Value = "bacbd3"

# --- Row 97 ---------------------------------------------------------------
$ws.Range("A97").Value = "Mod A - Fixes energy as driver data"
$ws.Range("B97").Value = 93
$ws.Range("C97").Value = "Rachel Hoesly"
$ws.Range("D97").Value = "Committed"
$ws.Range("E97").Value = 42346
$ws.Range("E97").NumberFormat = $dateFmt
$ws.Range("G97").Value = 42346
$ws.Range("G97").NumberFormat = $dateFmt
$ws.Range("H97").Value = "0f6b74f"

# --- View state: scroll the frozen pane down to the new rows and move the
#     selection to match where the author's cursor ended up after typing.
$win = $excel.ActiveWindow
$win.ScrollRow = 92
$win.ScrollColumn = 1
$ws.Range("G100").Select() | Out-Null
